$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 283
$ws.Range("I12").Value = 283
$ws.Range("K12").Value = 283
$ws.Range("M12").Value = -113
$ws.Range("H32").Value = 6343.769
$ws.Range("I32").Value = 6256.3335
$ws.Range("J32").Value = 6418.7144
$ws.Range("K32").Value = 6256.3335
$ws.Range("L32").Value = 6418.7144
$ws.Range("M32").Value = -5930.3335
$ws.Range("N32").Value = -7070.7144
$ws.Range("H34").Value = 8133.3335
$ws.Range("I34").Value = 8133.3335
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 8133.3335
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = $null
$ws.Range("N34").Value = $null
$ws.Range("H36").Value = 8133.3335
$ws.Range("I36").Value = 8133.3335
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 8133.3335
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = $null
$ws.Range("N36").Value = $null
$ws.Range("H55").Value = 457.37143
$ws.Range("I55").Value = 387
$ws.Range("K55").Value = 387
$ws.Range("M55").Value = -173
$ws.Range("H62").Value = 6819.1113
$ws.Range("I62").Value = 4988.9
$ws.Range("K62").Value = 4988.9
$ws.Range("M62").Value = -4364.9
$ws.Range("H65").Value = 6819.1113
$ws.Range("I65").Value = 4988.9
$ws.Range("K65").Value = 24944.5
$ws.Range("M65").Value = -21824.5
$ws.Range("H69").Value = 7659.619
$ws.Range("I69").Value = 7341.6
$ws.Range("J69").Value = 7759
$ws.Range("K69").Value = 22024.8
$ws.Range("L69").Value = 23277
$ws.Range("M69").Value = -21150.8
$ws.Range("N69").Value = -25025
$ws.Range("H72").Value = 7659.619
$ws.Range("I72").Value = 7341.6
$ws.Range("J72").Value = 7759
$ws.Range("K72").Value = 66074.40000000001
$ws.Range("L72").Value = 69831
$ws.Range("M72").Value = -61706.40000000001
$ws.Range("N72").Value = -78567
$ws.Range("H92").Value = 6230.25
$ws.Range("J92").Value = 6901.4287
$ws.Range("L92").Value = 6901.4287
$ws.Range("N92").Value = -9397.4287
$ws.Range("H98").Value = 557287.25
$ws.Range("I98").Value = 3349.1667
$ws.Range("K98").Value = 3349.1667
$ws.Range("M98").Value = -1851.1667
$ws.Range("H106").Value = 6556.3335
$ws.Range("I106").Value = 7000
$ws.Range("J106").Value = 6500.875
$ws.Range("K106").Value = 7000
$ws.Range("L106").Value = 6500.875
$ws.Range("M106").Value = -6369
$ws.Range("N106").Value = -7762.875
$ws.Range("H115").Value = 545
$ws.Range("I115").Value = 254
$ws.Range("J115").Value = 2000
$ws.Range("K115").Value = 762
$ws.Range("L115").Value = 6000
$ws.Range("M115").Value = 805
$ws.Range("N115").Value = -9134
$ws.Range("H122").Value = 557287.25
$ws.Range("I122").Value = 3349.1667
$ws.Range("K122").Value = 10047.5001
$ws.Range("M122").Value = -7597.500100000001
$ws.Range("H125").Value = 3999.2856
$ws.Range("J125").Value = 4799
$ws.Range("L125").Value = 43191
$ws.Range("N125").Value = -48111
$ws.Range("H129").Value = 981.7143
$ws.Range("I129").Value = 467.63635
$ws.Range("K129").Value = 1402.90905
$ws.Range("M129").Value = 3597.09095
$ws.Range("H132").Value = 1336.9259
$ws.Range("I132").Value = 1161.0435
$ws.Range("J132").Value = 2348.25
$ws.Range("K132").Value = 3483.1305
$ws.Range("L132").Value = 7044.75
$ws.Range("M132").Value = -953.1305000000002
$ws.Range("N132").Value = -12104.75
$ws.Range("H135").Value = 1103
$ws.Range("I135").Value = 832.4286
$ws.Range("J135").Value = 2997
$ws.Range("K135").Value = 7491.8574
$ws.Range("L135").Value = 26973
$ws.Range("M135").Value = -4956.8574
$ws.Range("N135").Value = -32043
$ws.Range("H137").Value = 3011.7031
$ws.Range("I137").Value = 2145.5
$ws.Range("J137").Value = 3350.652
$ws.Range("K137").Value = 6436.5
$ws.Range("L137").Value = 10051.956
$ws.Range("M137").Value = -3886.5
$ws.Range("N137").Value = -15151.956
$ws.Range("H138").Value = 2499.6516
$ws.Range("I138").Value = 1068.2354
$ws.Range("J138").Value = 3384.5273
$ws.Range("K138").Value = 3204.7062
$ws.Range("L138").Value = 10153.5819
$ws.Range("M138").Value = 1935.2938
$ws.Range("N138").Value = -20433.5819
$ws.Range("H141").Value = 2596.818
$ws.Range("I141").Value = 2771.25
$ws.Range("J141").Value = 2131.6667
$ws.Range("K141").Value = 8313.75
$ws.Range("L141").Value = 6395.000100000001
$ws.Range("M141").Value = -3133.75
$ws.Range("N141").Value = -16755.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4871.3945
$ws.Range("I32").Value = 3746.4
$ws.Range("K32").Value = 3746.4
$ws.Range("M32").Value = -3459.4
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = $null
$ws.Range("N42").Value = $null
$ws.Range("H61").Value = 5787.2173
$ws.Range("I61").Value = 4731.4546
$ws.Range("K61").Value = 4731.4546
$ws.Range("M61").Value = -4519.4546
$ws.Range("H63").Value = 5486.6113
$ws.Range("I63").Value = 3412.25
$ws.Range("J63").Value = 7146.1
$ws.Range("K63").Value = 3412.25
$ws.Range("L63").Value = 7146.1
$ws.Range("M63").Value = -2726.25
$ws.Range("N63").Value = -8518.1
$ws.Range("H66").Value = 5486.6113
$ws.Range("I66").Value = 3412.25
$ws.Range("J66").Value = 7146.1
$ws.Range("K66").Value = 17061.25
$ws.Range("L66").Value = 35730.5
$ws.Range("M66").Value = -13629.25
$ws.Range("N66").Value = -42594.5
$ws.Range("H74").Value = 41669380
$ws.Range("I74").Value = 41669380
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 41669380
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = $null
$ws.Range("N74").Value = $null
$ws.Range("H77").Value = 41669380
$ws.Range("I77").Value = 41669380
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 208346900
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = $null
$ws.Range("N77").Value = $null
$ws.Range("H102").Value = 2622
$ws.Range("I102").Value = 1699.6666
$ws.Range("J102").Value = 4005.5
$ws.Range("K102").Value = 1699.6666
$ws.Range("L102").Value = 4005.5
$ws.Range("M102").Value = -77.66660000000002
$ws.Range("N102").Value = -7249.5
$ws.Range("H110").Value = 4024.0435
$ws.Range("I110").Value = 3359.0476
$ws.Range("J110").Value = 11006.5
$ws.Range("K110").Value = 3359.0476
$ws.Range("L110").Value = 11006.5
$ws.Range("M110").Value = -1314.0476
$ws.Range("N110").Value = -15096.5
$ws.Range("H128").Value = 77500
$ws.Range("J128").Value = 77500
$ws.Range("L128").Value = 77500
$ws.Range("N128").Value = -87460
$ws.Range("H132").Value = 2924.3704
$ws.Range("I132").Value = 2081.2083
$ws.Range("K132").Value = 6243.624899999999
$ws.Range("M132").Value = -3713.624899999999
$ws.Range("H136").Value = 5787.2173
$ws.Range("I136").Value = 4731.4546
$ws.Range("K136").Value = 14194.3638
$ws.Range("M136").Value = -11644.3638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 21600.8
$ws.Range("I8").Value = 25001
$ws.Range("J8").Value = 8000
$ws.Range("K8").Value = 25001
$ws.Range("L8").Value = 8000
$ws.Range("M8").Value = -24861
$ws.Range("N8").Value = -8280
$ws.Range("H25").Value = 11500
$ws.Range("I25").Value = 5000
$ws.Range("K25").Value = 5000
$ws.Range("M25").Value = -4765
$ws.Range("H38").Value = 10990
$ws.Range("I38").Value = 10990
$ws.Range("K38").Value = 10990
$ws.Range("M38").Value = -10574
$ws.Range("H87").Value = 94750
$ws.Range("J87").Value = 94750
$ws.Range("L87").Value = 94750
$ws.Range("N87").Value = -97246
$ws.Range("H90").Value = 94750
$ws.Range("J90").Value = 94750
$ws.Range("L90").Value = 284250
$ws.Range("N90").Value = -296730
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = $null
$ws.Range("N92").Value = $null
$ws.Range("H99").Value = 2900.2856
$ws.Range("J99").Value = 2999.3333
$ws.Range("L99").Value = 2999.3333
$ws.Range("N99").Value = -5995.3333
$ws.Range("H105").Value = 21713
$ws.Range("I105").Value = 21648.4
$ws.Range("J105").Value = 21874.5
$ws.Range("K105").Value = 21648.4
$ws.Range("L105").Value = 21874.5
$ws.Range("M105").Value = -19901.4
$ws.Range("N105").Value = -25368.5
$ws.Range("H106").Value = 60559.75
$ws.Range("J106").Value = 60559.75
$ws.Range("L106").Value = 60559.75
$ws.Range("N106").Value = -63083.75
$ws.Range("H134").Value = 5557
$ws.Range("I134").Value = 4971.364
$ws.Range("J134").Value = 11999
$ws.Range("K134").Value = 14914.092
$ws.Range("L134").Value = 35997
$ws.Range("M134").Value = -12379.092
$ws.Range("N134").Value = -41067

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1727.3334
$ws.Range("I19").Value = 1011.5
$ws.Range("J19").Value = 2300
$ws.Range("K19").Value = 1011.5
$ws.Range("L19").Value = 2300
$ws.Range("M19").Value = -841.5
$ws.Range("N19").Value = -2640
$ws.Range("H24").Value = 1727.3334
$ws.Range("I24").Value = 1011.5
$ws.Range("J24").Value = 2300
$ws.Range("K24").Value = 1011.5
$ws.Range("L24").Value = 2300
$ws.Range("M24").Value = -841.5
$ws.Range("N24").Value = -2640
$ws.Range("H31").Value = 50884.047
$ws.Range("I31").Value = 1738.5454
$ws.Range("J31").Value = 100029.55
$ws.Range("K31").Value = 1738.5454
$ws.Range("L31").Value = 100029.55
$ws.Range("M31").Value = -1443.5454
$ws.Range("N31").Value = -100619.55
$ws.Range("H34").Value = 50884.047
$ws.Range("I34").Value = 1738.5454
$ws.Range("J34").Value = 100029.55
$ws.Range("K34").Value = 1738.5454
$ws.Range("L34").Value = 100029.55
$ws.Range("M34").Value = -1536.5454
$ws.Range("N34").Value = -100433.55
$ws.Range("H58").Value = 10758.714
$ws.Range("I58").Value = 2324.2
$ws.Range("K58").Value = 2324.2
$ws.Range("M58").Value = -2121.2
$ws.Range("H99").Value = 4699.5
$ws.Range("I99").Value = 4566
$ws.Range("K99").Value = 4566
$ws.Range("M99").Value = -3068
$ws.Range("H100").Value = 65000
$ws.Range("J100").Value = 65000
$ws.Range("L100").Value = 65000
$ws.Range("N100").Value = -67164
$ws.Range("H107").Value = 2516.4285
$ws.Range("I107").Value = 1304.3334
$ws.Range("J107").Value = 4698.2
$ws.Range("K107").Value = 1304.3334
$ws.Range("L107").Value = 4698.2
$ws.Range("M107").Value = 615.6666
$ws.Range("N107").Value = -8538.200000000001
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = $null
$ws.Range("N119").Value = $null
$ws.Range("H126").Value = 4699.5
$ws.Range("I126").Value = 4566
$ws.Range("K126").Value = 13698
$ws.Range("M126").Value = -11228
$ws.Range("H136").Value = 10758.714
$ws.Range("I136").Value = 2324.2
$ws.Range("K136").Value = 6972.599999999999
$ws.Range("M136").Value = -4422.599999999999
$ws.Range("H141").Value = 300872.38
$ws.Range("J141").Value = 300872.38
$ws.Range("L141").Value = 300872.38
$ws.Range("N141").Value = -311232.38

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 100245
$ws.Range("J2").Value = 100245
$ws.Range("L2").Value = 601470
$ws.Range("N2").Value = -601696
$ws.Range("H4").Value = 6026039
$ws.Range("J4").Value = 7658731.5
$ws.Range("L4").Value = 22976194.5
$ws.Range("N4").Value = -22976418.5
$ws.Range("H5").Value = 4565.8
$ws.Range("I5").Value = 274.66666
$ws.Range("J5").Value = 11002.5
$ws.Range("K5").Value = 823.9999799999999
$ws.Range("L5").Value = 33007.5
$ws.Range("M5").Value = -711.9999799999999
$ws.Range("N5").Value = -33231.5
$ws.Range("H14").Value = 4999.75
$ws.Range("I14").Value = 4999.75
$ws.Range("K14").Value = 14999.25
$ws.Range("M14").Value = -14826.25
$ws.Range("H38").Value = 47.615383
$ws.Range("I38").Value = 58.285713
$ws.Range("K38").Value = 174.857139
$ws.Range("M38").Value = 172.142861
$ws.Range("H62").Value = 833
$ws.Range("I62").Value = 833
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2499
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = $null
$ws.Range("H65").Value = 833
$ws.Range("I65").Value = 833
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 7497
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = $null
$ws.Range("H68").Value = 1699
$ws.Range("J68").Value = 2023.75
$ws.Range("L68").Value = 6071.25
$ws.Range("N68").Value = -7693.25
$ws.Range("H71").Value = 1699
$ws.Range("J71").Value = 2023.75
$ws.Range("L71").Value = 18213.75
$ws.Range("N71").Value = -26325.75
$ws.Range("H82").Value = 4999.5
$ws.Range("I82").Value = 4999
$ws.Range("K82").Value = 14997
$ws.Range("M82").Value = -14591
$ws.Range("H85").Value = 4999.5
$ws.Range("I85").Value = 4999
$ws.Range("K85").Value = 14997
$ws.Range("M85").Value = -13593
$ws.Range("H107").Value = 1007.0714
$ws.Range("I107").Value = 466.66666
$ws.Range("J107").Value = 1979.8
$ws.Range("K107").Value = 1399.99998
$ws.Range("L107").Value = 5939.4
$ws.Range("M107").Value = 520.0000199999999
$ws.Range("N107").Value = -9779.4
$ws.Range("H132").Value = 4916.905
$ws.Range("I132").Value = 4009.5557
$ws.Range("J132").Value = 5597.4165
$ws.Range("K132").Value = 36086.0013
$ws.Range("L132").Value = 50376.7485
$ws.Range("M132").Value = -33556.0013
$ws.Range("N132").Value = -55436.7485
$ws.Range("H134").Value = 12359.5
$ws.Range("I134").Value = 15131.75
$ws.Range("K134").Value = 45395.25
$ws.Range("M134").Value = -40325.25
$ws.Range("H135").Value = 4565.8
$ws.Range("I135").Value = 274.66666
$ws.Range("J135").Value = 11002.5
$ws.Range("K135").Value = 2471.99994
$ws.Range("L135").Value = 99022.5
$ws.Range("M135").Value = 63.0000600000003
$ws.Range("N135").Value = -104092.5
$ws.Range("H136").Value = 3309.6667
$ws.Range("I136").Value = 2464.5
$ws.Range("K136").Value = 7393.5
$ws.Range("M136").Value = -2293.5
$ws.Range("H137").Value = 67993.60000000001
$ws.Range("I137").Value = 718.8182
$ws.Range("J137").Value = 252999.25
$ws.Range("K137").Value = 2156.4546
$ws.Range("L137").Value = 758997.75
$ws.Range("M137").Value = 2943.5454
$ws.Range("N137").Value = -769197.75
$ws.Range("H138").Value = 1923.3334
$ws.Range("I138").Value = 1923.3334
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 5770.0002
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = $null
$ws.Range("N138").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 19999
$ws.Range("J38").Value = 19999
$ws.Range("L38").Value = 19999
$ws.Range("N38").Value = -20925
$ws.Range("H40").Value = 24666.666
$ws.Range("J40").Value = 30000
$ws.Range("L40").Value = 30000
$ws.Range("N40").Value = -30302
$ws.Range("H62").Value = 45506
$ws.Range("J62").Value = 49998.668
$ws.Range("L62").Value = 49998.668
$ws.Range("N62").Value = -51370.668
$ws.Range("H65").Value = 45506
$ws.Range("J65").Value = 49998.668
$ws.Range("L65").Value = 149996.004
$ws.Range("N65").Value = -156860.004
$ws.Range("H102").Value = 2939.4614
$ws.Range("I102").Value = 2339.5625
$ws.Range("K102").Value = 2339.5625
$ws.Range("M102").Value = -717.5625
$ws.Range("H107").Value = 660.94446
$ws.Range("J107").Value = 463.75
$ws.Range("L107").Value = 463.75
$ws.Range("N107").Value = -4303.75
$ws.Range("H122").Value = 8116.5586
$ws.Range("I122").Value = 10418.375
$ws.Range("K122").Value = 31255.125
$ws.Range("M122").Value = -28805.125
$ws.Range("H126").Value = 7001.7144
$ws.Range("I126").Value = 4998
$ws.Range("K126").Value = 14994
$ws.Range("M126").Value = -12524
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = $null
$ws.Range("N128").Value = $null
$ws.Range("H132").Value = 2358.56
$ws.Range("I132").Value = 1307.625
$ws.Range("K132").Value = 3922.875
$ws.Range("M132").Value = -1392.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10249.833
$ws.Range("I40").Value = 6250
$ws.Range("J40").Value = 12249.75
$ws.Range("K40").Value = 6250
$ws.Range("L40").Value = 12249.75
$ws.Range("M40").Value = -6114
$ws.Range("N40").Value = -12521.75
$ws.Range("H46").Value = 4534.722
$ws.Range("J46").Value = 5179.591
$ws.Range("L46").Value = 5179.591
$ws.Range("N46").Value = -5555.591
$ws.Range("H68").Value = 5612.25
$ws.Range("I68").Value = 4599.5
$ws.Range("J68").Value = 6625
$ws.Range("K68").Value = 4599.5
$ws.Range("L68").Value = 6625
$ws.Range("M68").Value = -3850.5
$ws.Range("N68").Value = -8123
$ws.Range("H71").Value = 5612.25
$ws.Range("I71").Value = 4599.5
$ws.Range("J71").Value = 6625
$ws.Range("K71").Value = 22997.5
$ws.Range("L71").Value = 33125
$ws.Range("M71").Value = -19253.5
$ws.Range("N71").Value = -40613
$ws.Range("H130").Value = 84900
$ws.Range("J130").Value = 84900
$ws.Range("L130").Value = 84900
$ws.Range("N130").Value = -94940
$ws.Range("H136").Value = 10880.059
$ws.Range("I136").Value = 8457.77
$ws.Range("K136").Value = 25373.31
$ws.Range("M136").Value = -22823.31

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 22501.75
$ws.Range("J18").Value = 22501.75
$ws.Range("L18").Value = 22501.75
$ws.Range("N18").Value = -22847.75
$ws.Range("H24").Value = 10000
$ws.Range("J24").Value = 10000
$ws.Range("L24").Value = 10000
$ws.Range("N24").Value = -10460
$ws.Range("H75").Value = 47500
$ws.Range("J75").Value = 75000
$ws.Range("L75").Value = 75000
$ws.Range("N75").Value = -76872
$ws.Range("H78").Value = 47500
$ws.Range("J78").Value = 75000
$ws.Range("L78").Value = 225000
$ws.Range("N78").Value = -234360
$ws.Range("H81").Value = 5383.3125
$ws.Range("I81").Value = 4116.5
$ws.Range("J81").Value = 5564.2856
$ws.Range("K81").Value = 8233
$ws.Range("L81").Value = 11128.5712
$ws.Range("M81").Value = -7172
$ws.Range("N81").Value = -13250.5712
$ws.Range("H84").Value = 5383.3125
$ws.Range("I84").Value = 4116.5
$ws.Range("J84").Value = 5564.2856
$ws.Range("K84").Value = 41165
$ws.Range("L84").Value = 55642.856
$ws.Range("M84").Value = -35861
$ws.Range("N84").Value = -66250.856
$ws.Range("H100").Value = 883.5833
$ws.Range("I100").Value = 782.0909
$ws.Range("K100").Value = 1564.1818
$ws.Range("M100").Value = -1023.1818
$ws.Range("H122").Value = 3213.1904
$ws.Range("I122").Value = 1863
$ws.Range("J122").Value = 4225.8335
$ws.Range("K122").Value = 5589
$ws.Range("L122").Value = 12677.5005
$ws.Range("M122").Value = -3139
$ws.Range("N122").Value = -17577.5005
$ws.Range("H126").Value = 1665.9375
$ws.Range("I126").Value = 1776.5
$ws.Range("K126").Value = 5329.5
$ws.Range("M126").Value = -2859.5
$ws.Range("H132").Value = 1757.5714
$ws.Range("I132").Value = 1037.75
$ws.Range("J132").Value = 6076.5
$ws.Range("K132").Value = 3113.25
$ws.Range("L132").Value = 18229.5
$ws.Range("M132").Value = -583.25
$ws.Range("N132").Value = -23289.5
$ws.Range("H136").Value = 3624.0833
$ws.Range("J136").Value = 5664.222
$ws.Range("L136").Value = 16992.666
$ws.Range("N136").Value = -22092.666
